# Apply updated Back/Lay odds values to rows 2-7 of the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.66
$ws.Range("G2").Value = 1.77
$ws.Range("H2").Value = 5.8
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 3.6
$ws.Range("K2").Value = 4.1
$ws.Range("L2").Value = 1.45
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 3.4
$ws.Range("O2").Value = 1.35
$ws.Range("P2").Value = 1.81
$ws.Range("Q2").Value = 2.04
$ws.Range("R2").Value = 1.31
$ws.Range("S2").Value = 3.85
$ws.Range("T2").Value = 1.98
$ws.Range("U2").Value = 1.83
$ws.Range("V2").Value = 1.17
$ws.Range("W2").Value = 2.28
$ws.Range("X2").Value = 13
$ws.Range("Z2").Value = 50
$ws.Range("AB2").Value = 7.8
$ws.Range("AC2").Value = 8.800000000000001
$ws.Range("AD2").Value = 26
$ws.Range("AE2").Value = 110
$ws.Range("AF2").Value = 10
$ws.Range("AH2").Value = 22
$ws.Range("AI2").Value = 120
$ws.Range("AJ2").Value = 17.5
$ws.Range("AK2").Value = 20
$ws.Range("AL2").Value = 44
$ws.Range("AO2").Value = 150

# Row 3
$ws.Range("F3").Value = 1.46
$ws.Range("G3").Value = 1.5
$ws.Range("H3").Value = 7.8
$ws.Range("I3").Value = 9.4
$ws.Range("J3").Value = 4.7
$ws.Range("K3").Value = 5.1
$ws.Range("L3").Value = 1.36
$ws.Range("O3").Value = 1.26
$ws.Range("P3").Value = 2.18
$ws.Range("Q3").Value = 1.75
$ws.Range("R3").Value = 1.45
$ws.Range("S3").Value = 2.96
$ws.Range("T3").Value = 1.91
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.12
$ws.Range("W3").Value = 2.96
$ws.Range("X3").Value = 20
$ws.Range("Y3").Value = 30
$ws.Range("Z3").Value = 75
$ws.Range("AA3").Value = 340
$ws.Range("AB3").Value = 9.199999999999999
$ws.Range("AC3").Value = 11
$ws.Range("AD3").Value = 32
$ws.Range("AE3").Value = 130
$ws.Range("AF3").Value = 9
$ws.Range("AH3").Value = 24
$ws.Range("AI3").Value = 120
$ws.Range("AJ3").Value = 13
$ws.Range("AK3").Value = 15
$ws.Range("AL3").Value = 36
$ws.Range("AM3").Value = 150
$ws.Range("AN3").Value = 7.4
$ws.Range("AO3").Value = 160

# Row 4
$ws.Range("F4").Value = 2.08
$ws.Range("G4").Value = 2.14
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 3.9
$ws.Range("N4").Value = 4.7
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 2.28
$ws.Range("Q4").Value = 1.69
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 2.74
$ws.Range("U4").Value = 2.34
$ws.Range("V4").Value = 1.35
$ws.Range("W4").Value = 1.87
$ws.Range("Y4").Value = 18
$ws.Range("Z4").Value = 30
$ws.Range("AC4").Value = 9
$ws.Range("AD4").Value = 16
$ws.Range("AE4").Value = 40
$ws.Range("AF4").Value = 14.5
$ws.Range("AG4").Value = 10.5
$ws.Range("AH4").Value = 16.5
$ws.Range("AI4").Value = 120
$ws.Range("AJ4").Value = 26
$ws.Range("AK4").Value = 20
$ws.Range("AL4").Value = 30
$ws.Range("AM4").Value = 200
$ws.Range("AN4").Value = 12.5
$ws.Range("AO4").Value = 32

# Row 5
$ws.Range("F5").Value = 2.32
$ws.Range("G5").Value = 2.5
$ws.Range("H5").Value = 3.55
$ws.Range("K5").Value = 3.35
$ws.Range("N5").Value = 2.74
$ws.Range("O5").Value = 1.48
$ws.Range("P5").Value = 1.59
$ws.Range("Q5").Value = 2.44
$ws.Range("S5").Value = 5
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 1.33
$ws.Range("W5").Value = 1.67
$ws.Range("X5").Value = 10
$ws.Range("Y5").Value = 980
$ws.Range("AB5").Value = 8.199999999999999
$ws.Range("AC5").Value = 7.6
$ws.Range("AF5").Value = 26
$ws.Range("AH5").Value = 23
$ws.Range("AI5").Value = 90
$ws.Range("AK5").Value = 32
$ws.Range("AM5").Value = 200
$ws.Range("AN5").Value = 32

# Row 6
$ws.Range("F6").Value = 1.94
$ws.Range("G6").Value = 2.14
$ws.Range("H6").Value = 4.8
$ws.Range("I6").Value = 5.5
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 3.4
$ws.Range("L6").Value = 1.61
$ws.Range("N6").Value = 2.46
$ws.Range("P6").Value = 1.48
$ws.Range("Q6").Value = 2.72
$ws.Range("R6").Value = 1.16
$ws.Range("S6").Value = 5.7
$ws.Range("T6").Value = 2.22
$ws.Range("U6").Value = 1.65
$ws.Range("V6").Value = 1.22
$ws.Range("W6").Value = 1.88
$ws.Range("X6").Value = 90
$ws.Range("AB6").Value = 29
$ws.Range("AC6").Value = 42
$ws.Range("AF6").Value = 40
$ws.Range("AJ6").Value = 220
$ws.Range("AK6").Value = 140

# Row 7
$ws.Range("F7").Value = 3.6
$ws.Range("I7").Value = 2.4
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 2.8
